# tv_generation_catalog.xlsx — "adding comment headers to all py files"
#
# On Sheet1, column B ("Make") flips from 1 to 0 for every data row except
# the very first one (row 2 stays at 1); rows 3-13 go from 1 -> 0.
# The active selection also moves from D20 to C16.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Flip the "Make" flag (column B) to 0 for rows 3 through 13.
$ws.Range("B3:B13").Value = 0

# Update the sheet's saved selection/active cell.
$ws.Range("C16").Select()
